# "Melhora na função de requisição"
# The product-listing sheet is regenerated from the (improved) API request:
# the two previous rows of product data (7 columns each) are replaced by a
# simple single-column list of product codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe out the old data rows (everything below the header row, across all
# previously-used columns A:G) before writing the new result set.
$ws.Range("A2:G3").ClearContents() | Out-Null

# New result set returned by the request: just a code per row, in column A.
$codes = @("C-2184", "AC 30937", "HG 30784", "BD3442", "BD4190")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

# Leave the selection back at the top of the sheet.
$ws.Range("A1").Select() | Out-Null
